# Scheduled runner update: refresh market-price derived columns (H:N) for
# a set of Leve rows across multiple job sheets. Values below were produced
# by the upstream price-scraping job; this script just writes them in place.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3505.1333
$ws.Range("I64").Value = 2982.4285
$ws.Range("K64").Value = 2982.4285
$ws.Range("M64").Value = -2734.4285
$ws.Range("H67").Value = 3505.1333
$ws.Range("I67").Value = 2982.4285
$ws.Range("K67").Value = 2982.4285
$ws.Range("M67").Value = -2124.4285
$ws.Range("H111").Value = 4024.9285
$ws.Range("I111").Value = 3581.25
$ws.Range("J111").Value = 4616.5
$ws.Range("K111").Value = 10743.75
$ws.Range("L111").Value = 13849.5
$ws.Range("M111").Value = -7676.75
$ws.Range("N111").Value = -19983.5
$ws.Range("H125").Value = 1140
$ws.Range("J125").Value = 1140
$ws.Range("L125").Value = 10260
$ws.Range("N125").Value = -15180
$ws.Range("H132").Value = 2529.9055
$ws.Range("I132").Value = 1240.3385
$ws.Range("J132").Value = 11843.444
$ws.Range("K132").Value = 3721.0155
$ws.Range("L132").Value = 35530.33199999999
$ws.Range("M132").Value = -1191.0155
$ws.Range("N132").Value = -40590.33199999999
$ws.Range("H135").Value = 27778576
$ws.Range("I135").Value = 567.76
$ws.Range("J135").Value = 90910420
$ws.Range("K135").Value = 5109.84
$ws.Range("L135").Value = 818193780
$ws.Range("M135").Value = -2574.84
$ws.Range("N135").Value = -818198850
$ws.Range("H138").Value = 2249496.2
$ws.Range("I138").Value = 1051.625
$ws.Range("J138").Value = 3511780.8
$ws.Range("K138").Value = 3154.875
$ws.Range("L138").Value = 10535342.4
$ws.Range("M138").Value = 1985.125
$ws.Range("N138").Value = -10545622.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3200.0557
$ws.Range("I63").Value = 2800.1428
$ws.Range("K63").Value = 2800.1428
$ws.Range("M63").Value = -2114.1428
$ws.Range("H66").Value = 3200.0557
$ws.Range("I66").Value = 2800.1428
$ws.Range("K66").Value = 14000.714
$ws.Range("M66").Value = -10568.714
$ws.Range("H97").Value = 744.06665
$ws.Range("I97").Value = 687.4
$ws.Range("J97").Value = 857.4
$ws.Range("K97").Value = 687.4
$ws.Range("L97").Value = 857.4
$ws.Range("M97").Value = -191.4
$ws.Range("N97").Value = -1849.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 27000
$ws.Range("J40").Value = 27000
$ws.Range("L40").Value = 27000
$ws.Range("N40").Value = -27530
$ws.Range("H94").Value = 1210.9333
$ws.Range("I94").Value = 1012.8333
$ws.Range("J94").Value = 2003.3334
$ws.Range("K94").Value = 1012.8333
$ws.Range("L94").Value = 2003.3334
$ws.Range("M94").Value = -561.8333
$ws.Range("N94").Value = -2905.3334
$ws.Range("H96").Value = 16250
$ws.Range("I96").Value = 2500
$ws.Range("J96").Value = 30000
$ws.Range("K96").Value = 2500
$ws.Range("L96").Value = 30000
$ws.Range("M96").Value = 246
$ws.Range("N96").Value = -35492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 970.8
$ws.Range("I16").Value = 766
$ws.Range("J16").Value = 1058.5714
$ws.Range("K16").Value = 766
$ws.Range("L16").Value = 1058.5714
$ws.Range("M16").Value = -479
$ws.Range("N16").Value = -1632.5714
$ws.Range("H31").Value = 1384.2063
$ws.Range("I31").Value = 862.7
$ws.Range("J31").Value = 2291.1738
$ws.Range("K31").Value = 862.7
$ws.Range("L31").Value = 2291.1738
$ws.Range("M31").Value = -567.7
$ws.Range("N31").Value = -2881.1738
$ws.Range("H34").Value = 1384.2063
$ws.Range("I34").Value = 862.7
$ws.Range("J34").Value = 2291.1738
$ws.Range("K34").Value = 862.7
$ws.Range("L34").Value = 2291.1738
$ws.Range("M34").Value = -660.7
$ws.Range("N34").Value = -2695.1738
$ws.Range("H62").Value = 2355.5557
$ws.Range("I62").Value = 2300
$ws.Range("K62").Value = 2300
$ws.Range("M62").Value = -1676
$ws.Range("H65").Value = 2355.5557
$ws.Range("I65").Value = 2300
$ws.Range("K65").Value = 11500
$ws.Range("M65").Value = -8380
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H99").Value = 3917.6765
$ws.Range("I99").Value = 3669.84
$ws.Range("J99").Value = 4606.1113
$ws.Range("K99").Value = 3669.84
$ws.Range("L99").Value = 4606.1113
$ws.Range("M99").Value = -2171.84
$ws.Range("N99").Value = -7602.1113
$ws.Range("H113").Value = 970.8
$ws.Range("I113").Value = 766
$ws.Range("J113").Value = 1058.5714
$ws.Range("K113").Value = 766
$ws.Range("L113").Value = 1058.5714
$ws.Range("M113").Value = 1404
$ws.Range("N113").Value = -5398.5714
$ws.Range("H126").Value = 3917.6765
$ws.Range("I126").Value = 3669.84
$ws.Range("J126").Value = 4606.1113
$ws.Range("K126").Value = 11009.52
$ws.Range("L126").Value = 13818.3339
$ws.Range("M126").Value = -8539.52
$ws.Range("N126").Value = -18758.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4942.5713
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4942.5713
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 14827.7139
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16199.7139
$ws.Range("H65").Value = 4942.5713
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4942.5713
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 44483.14169999999
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -51347.14169999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H97").Value = 683.2083
$ws.Range("I97").Value = 666.5238000000001
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 666.5238000000001
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -170.5238000000001
$ws.Range("N97").Value = -1792
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178
$ws.Range("H132").Value = 1697032.8
$ws.Range("I132").Value = 2005.4565
$ws.Range("J132").Value = 7694821.5
$ws.Range("K132").Value = 6016.3695
$ws.Range("L132").Value = 23084464.5
$ws.Range("M132").Value = -3486.3695
$ws.Range("N132").Value = -23089524.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 29194.5
$ws.Range("J92").Value = 29194.5
$ws.Range("L92").Value = 29194.5
$ws.Range("N92").Value = -34186.5
$ws.Range("H96").Value = 28375
$ws.Range("J96").Value = 28375
$ws.Range("L96").Value = 28375
$ws.Range("N96").Value = -33867
$ws.Range("H136").Value = 1683.75
$ws.Range("I136").Value = 1230.909
$ws.Range("J136").Value = 2395.3572
$ws.Range("K136").Value = 3692.727
$ws.Range("L136").Value = 7186.071599999999
$ws.Range("M136").Value = -1142.727
$ws.Range("N136").Value = -12286.0716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 11500
$ws.Range("J54").Value = 11500
$ws.Range("L54").Value = 11500
$ws.Range("N54").Value = -12540
$ws.Range("H70").Value = 24665
$ws.Range("J70").Value = 25450
$ws.Range("L70").Value = 25450
$ws.Range("N70").Value = -26080
$ws.Range("H73").Value = 24665
$ws.Range("J73").Value = 25450
$ws.Range("L73").Value = 25450
$ws.Range("N73").Value = -27634
$ws.Range("H126").Value = 906.619
$ws.Range("I126").Value = 817.6923
$ws.Range("J126").Value = 1051.125
$ws.Range("K126").Value = 2453.0769
$ws.Range("L126").Value = 3153.375
$ws.Range("M126").Value = 16.92309999999998
$ws.Range("N126").Value = -8093.375
